$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.378.12"
$ws.Range("E2").Value = "  +1.50%  "

# Row 3
$ws.Range("D3").Value = "1.867.50"
$ws.Range("E3").Value = "  +1.45%  "

# Row 4
$ws.Range("E4").Value = "  +0.44%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.54"
$ws.Range("E5").Value = "  +3.39%  "

# Row 6
$ws.Range("E6").Value = "  +0.91%  "

# Row 7
$ws.Range("E7").Value = "  +0.31%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.78"
$ws.Range("E8").Value = "  +7.52%  "

# Row 9
$ws.Range("E9").Value = "  +0.67%  "

# Row 10
$ws.Range("E10").Value = "  +1.46%  "

# Row 11
$ws.Range("E11").Value = "  +0.76%  "

# Row 12
$ws.Range("D12").Value = "2.135.25"
$ws.Range("E12").Value = "  +1.32%  "

# Row 13
$ws.Range("D13").Value = "1.904.88"
$ws.Range("E13").Value = "  +3.48%  "

# Row 14
$ws.Range("E14").Value = "  +0.82%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.682"
$ws.Range("E15").Value = "  +1.31%  "

# Row 16
$ws.Range("E16").Value = "  +1.94%  "

# Row 17
$ws.Range("D17").Value = "35.353.08"
$ws.Range("E17").Value = "  +1.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.23"
$ws.Range("E18").Value = "  +0.50%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0801"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.85"
$ws.Range("E20").Value = "  +0.40%  "

# Row 21
$ws.Range("E21").Value = "  +0.81%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.77"
$ws.Range("E22").Value = "  +1.58%  "

# Row 23
$ws.Range("E23").Value = "  +0.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.24"
$ws.Range("E24").Value = "  -1.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.99"
$ws.Range("E25").Value = "  -0.85%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.27"
$ws.Range("E26").Value = "  +6.02%  "

# Row 27
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.90"
$ws.Range("E27").Value = "  +25.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.79"
$ws.Range("E28").Value = "  +1.83%  "

# Row 29
$ws.Range("E29").Value = "  +1.43%  "

# Row 30
$ws.Range("E30").Value = "  +2.02%  "

# Row 31
$ws.Range("E31").Value = "  +0.37%  "

# Row 32
$ws.Range("E32").Value = "  +1.95%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("E33").Value = "  +25.67%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.05"
$ws.Range("E34").Value = "  +2.35%  "

# Row 35
$ws.Range("E35").Value = "  +8.74%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.825"
$ws.Range("E36").Value = "  +18.70%  "

# Row 37
$ws.Range("E37").Value = "  +5.56%  "

# Row 38
$ws.Range("E38").Value = "  +3.42%  "

# Row 39
$ws.Range("E39").Value = "  +4.59%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "91.22"
$ws.Range("E40").Value = "  +1.00%  "

# Row 41
$ws.Range("D41").Value = "1.349.91"
$ws.Range("E41").Value = "  +0.27%  "

# Row 42
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.20"
$ws.Range("E42").Value = "  +1.71%  "

# Row 43
$ws.Range("B43").Value = "Gas"
$ws.Range("C43").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.36"
$ws.Range("E43").Value = "  +58.90%  "

# Row 44
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0600"
$ws.Range("E44").Value = "  +14.67%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.37"
$ws.Range("E45").Value = "  +2.83%  "

# Row 46
$ws.Range("E46").Value = "  +0.69%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.65"
$ws.Range("E47").Value = "  +6.27%  "

# Row 48
$ws.Range("E48").Value = "  -1.01%  "

# Row 49
$ws.Range("D49").Value = "2.049.90"
$ws.Range("E49").Value = "  +1.46%  "

# Row 50
$ws.Range("E50").Value = "  +3.32%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.42"
$ws.Range("E51").Value = "  -1.14%  "
